$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 189.4
$ws.Range("I33").Value = 212.23529
$ws.Range("K33").Value = 212.23529
$ws.Range("M33").Value = 16.76471000000001
$ws.Range("H86").Value = 5586.522
$ws.Range("I86").Value = 1206.4286
$ws.Range("J86").Value = 12400
$ws.Range("K86").Value = 1206.4286
$ws.Range("L86").Value = 12400
$ws.Range("M86").Value = -83.42859999999996
$ws.Range("N86").Value = -14646
$ws.Range("H88").Value = 550.2
$ws.Range("I88").Value = 357.5
$ws.Range("J88").Value = 678.6667
$ws.Range("K88").Value = 357.5
$ws.Range("L88").Value = 678.6667
$ws.Range("M88").Value = 48.5
$ws.Range("N88").Value = -1490.6667
$ws.Range("H89").Value = 5586.522
$ws.Range("I89").Value = 1206.4286
$ws.Range("J89").Value = 12400
$ws.Range("K89").Value = 6032.143
$ws.Range("L89").Value = 62000
$ws.Range("M89").Value = -416.143
$ws.Range("N89").Value = -73232
$ws.Range("H91").Value = 550.2
$ws.Range("I91").Value = 357.5
$ws.Range("J91").Value = 678.6667
$ws.Range("K91").Value = 357.5
$ws.Range("L91").Value = 678.6667
$ws.Range("M91").Value = 1046.5
$ws.Range("N91").Value = -3486.6667
$ws.Range("H100").Value = 1709.4546
$ws.Range("I100").Value = 934
$ws.Range("J100").Value = 2640
$ws.Range("K100").Value = 934
$ws.Range("L100").Value = 2640
$ws.Range("M100").Value = -393
$ws.Range("N100").Value = -3722
$ws.Range("H116").Value = 11908543
$ws.Range("I116").Value = 20834608
$ws.Range("J116").Value = 7122.778
$ws.Range("K116").Value = 20834608
$ws.Range("L116").Value = 7122.778
$ws.Range("M116").Value = -20831166
$ws.Range("N116").Value = -14006.778
$ws.Range("H129").Value = 726
$ws.Range("J129").Value = 805.4815
$ws.Range("L129").Value = 2416.4445
$ws.Range("N129").Value = -12416.4445
$ws.Range("H131").Value = 1997.04
$ws.Range("I131").Value = 1148
$ws.Range("J131").Value = 2780.7693
$ws.Range("K131").Value = 3444
$ws.Range("L131").Value = 8342.3079
$ws.Range("M131").Value = 1596
$ws.Range("N131").Value = -18422.3079

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1256.25
$ws.Range("I2").Value = 1208.3334
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 1208.3334
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = -1095.3334
$ws.Range("N2").Value = -1626
$ws.Range("H45").Value = 2520.04
$ws.Range("J45").Value = 2500.7273
$ws.Range("L45").Value = 2500.7273
$ws.Range("N45").Value = -3254.7273
$ws.Range("H63").Value = 2843174.5
$ws.Range("I63").Value = 2491
$ws.Range("J63").Value = 15626250
$ws.Range("K63").Value = 2491
$ws.Range("L63").Value = 15626250
$ws.Range("M63").Value = -1805
$ws.Range("N63").Value = -15627622
$ws.Range("H66").Value = 2843174.5
$ws.Range("I66").Value = 2491
$ws.Range("J66").Value = 15626250
$ws.Range("K66").Value = 12455
$ws.Range("L66").Value = 78131250
$ws.Range("M66").Value = -9023
$ws.Range("N66").Value = -78138114
$ws.Range("H97").Value = 1070.2
$ws.Range("I97").Value = 1125.25
$ws.Range("J97").Value = 850
$ws.Range("K97").Value = 1125.25
$ws.Range("L97").Value = 850
$ws.Range("M97").Value = -629.25
$ws.Range("N97").Value = -1842
$ws.Range("H116").Value = 1256.25
$ws.Range("I116").Value = 1208.3334
$ws.Range("J116").Value = 1400
$ws.Range("K116").Value = 1208.3334
$ws.Range("L116").Value = 1400
$ws.Range("M116").Value = 1085.6666
$ws.Range("N116").Value = -5988
$ws.Range("H132").Value = 13974.349
$ws.Range("I132").Value = 2295.7742
$ws.Range("J132").Value = 44144
$ws.Range("K132").Value = 6887.3226
$ws.Range("L132").Value = 132432
$ws.Range("M132").Value = -4357.3226
$ws.Range("N132").Value = -137492

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1256.25
$ws.Range("I3").Value = 1208.3334
$ws.Range("J3").Value = 1400
$ws.Range("K3").Value = 1208.3334
$ws.Range("L3").Value = 1400
$ws.Range("M3").Value = -1094.3334
$ws.Range("N3").Value = -1628
$ws.Range("H22").Value = 864.7857
$ws.Range("I22").Value = 877.46155
$ws.Range("K22").Value = 877.46155
$ws.Range("M22").Value = -704.46155
$ws.Range("H24").Value = 758
$ws.Range("I24").Value = 758
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 758
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -523
$ws.Range("N24").ClearContents()
$ws.Range("H94").Value = 828.8570999999999
$ws.Range("I94").Value = 694.3333
$ws.Range("J94").Value = 1165.1666
$ws.Range("K94").Value = 694.3333
$ws.Range("L94").Value = 1165.1666
$ws.Range("M94").Value = -243.3333
$ws.Range("N94").Value = -2067.1666
$ws.Range("H105").Value = 1065726.6
$ws.Range("I105").Value = 1584.2273
$ws.Range("J105").Value = 2002172
$ws.Range("K105").Value = 1584.2273
$ws.Range("L105").Value = 2002172
$ws.Range("M105").Value = 162.7727
$ws.Range("N105").Value = -2005666
$ws.Range("H134").Value = 3831.4412
$ws.Range("I134").Value = 3989.2188
$ws.Range("K134").Value = 11967.6564
$ws.Range("M134").Value = -9432.6564

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16107.323
$ws.Range("I58").Value = 1783.8572
$ws.Range("J58").Value = 26133.75
$ws.Range("K58").Value = 1783.8572
$ws.Range("L58").Value = 26133.75
$ws.Range("M58").Value = -1580.8572
$ws.Range("N58").Value = -26539.75
$ws.Range("H110").Value = 30702
$ws.Range("J110").Value = 30702
$ws.Range("L110").Value = 30702
$ws.Range("N110").Value = -38882
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H134").Value = 926.6923
$ws.Range("I134").Value = 857.6316
$ws.Range("K134").Value = 2572.8948
$ws.Range("M134").Value = -37.89480000000003
$ws.Range("H136").Value = 16107.323
$ws.Range("I136").Value = 1783.8572
$ws.Range("J136").Value = 26133.75
$ws.Range("K136").Value = 5351.571599999999
$ws.Range("L136").Value = 78401.25
$ws.Range("M136").Value = -2801.571599999999
$ws.Range("N136").Value = -83501.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 433.33334
$ws.Range("I7").Value = 150
$ws.Range("J7").Value = 575
$ws.Range("K7").Value = 450
$ws.Range("L7").Value = 1725
$ws.Range("M7").Value = -338
$ws.Range("N7").Value = -1949
$ws.Range("H122").Value = 945.8889
$ws.Range("I122").Value = 590
$ws.Range("K122").Value = 5310
$ws.Range("M122").Value = -2860
$ws.Range("H131").Value = 744.37
$ws.Range("J131").Value = 757.835
$ws.Range("L131").Value = 2273.505
$ws.Range("N131").Value = -12353.505

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8179.2856
$ws.Range("I113").Value = 11223.444
$ws.Range("J113").Value = 2699.8
$ws.Range("K113").Value = 11223.444
$ws.Range("L113").Value = 2699.8
$ws.Range("M113").Value = -9053.444
$ws.Range("N113").Value = -7039.8
$ws.Range("H122").Value = 4348.0527
$ws.Range("I122").Value = 4590.4443
$ws.Range("J122").Value = 4129.9
$ws.Range("K122").Value = 13771.3329
$ws.Range("L122").Value = 12389.7
$ws.Range("M122").Value = -11321.3329
$ws.Range("N122").Value = -17289.7
$ws.Range("H129").Value = 49694.25
$ws.Range("J129").Value = 49694.25
$ws.Range("L129").Value = 49694.25
$ws.Range("N129").Value = -59694.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 587.8
$ws.Range("I16").Value = 532.58826
$ws.Range("J16").Value = 900.6667
$ws.Range("K16").Value = 532.58826
$ws.Range("L16").Value = 900.6667
$ws.Range("M16").Value = -362.58826
$ws.Range("N16").Value = -1240.6667
$ws.Range("H22").Value = 5836.8
$ws.Range("I22").Value = 5811.6924
$ws.Range("J22").Value = 6000
$ws.Range("K22").Value = 5811.6924
$ws.Range("L22").Value = 6000
$ws.Range("M22").Value = -5516.6924
$ws.Range("N22").Value = -6590
$ws.Range("H27").Value = 5836.8
$ws.Range("I27").Value = 5811.6924
$ws.Range("J27").Value = 6000
$ws.Range("K27").Value = 5811.6924
$ws.Range("L27").Value = 6000
$ws.Range("M27").Value = -5704.6924
$ws.Range("N27").Value = -6214
$ws.Range("H40").Value = 3285.2917
$ws.Range("I40").Value = 3257.35
$ws.Range("J40").Value = 3425
$ws.Range("K40").Value = 3257.35
$ws.Range("L40").Value = 3425
$ws.Range("M40").Value = -3121.35
$ws.Range("N40").Value = -3697
$ws.Range("H68").Value = 2699.2
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 2999
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 2999
$ws.Range("M68").Value = -751
$ws.Range("N68").Value = -4497
$ws.Range("H71").Value = 2699.2
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 2999
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 14995
$ws.Range("M71").Value = -3756
$ws.Range("N71").Value = -22483
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H100").Value = 2077.4707
$ws.Range("I100").Value = 1196.4
$ws.Range("J100").Value = 2444.5833
$ws.Range("K100").Value = 1196.4
$ws.Range("L100").Value = 2444.5833
$ws.Range("M100").Value = -655.4000000000001
$ws.Range("N100").Value = -3526.5833
$ws.Range("H136").Value = 1677.3235
$ws.Range("I136").Value = 1622.8
$ws.Range("J136").Value = 2086.25
$ws.Range("K136").Value = 4868.4
$ws.Range("L136").Value = 6258.75
$ws.Range("M136").Value = -2318.4
$ws.Range("N136").Value = -11358.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 83333550
$ws.Range("I107").Value = 100000216
$ws.Range("J107").Value = 195
$ws.Range("K107").Value = 300000648
$ws.Range("L107").Value = 585
$ws.Range("M107").Value = -299998728
$ws.Range("N107").Value = -4425
$ws.Range("H136").Value = 28676992
$ws.Range("I136").Value = 38233252
$ws.Range("K136").Value = 114699756
$ws.Range("M136").Value = -114697206
